$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# B2: date text changes from 08/25/20 -> 09/04/20 (must remain plain text, not become a date serial)
$ws.Range("B2").Value = "'09/04/20"
$ws.Range("B2").Style = "Normal"

# D2: PRICE 3222 -> 12
$ws.Range("D2").Value = 12

# F2: NET_EFFECT_TO_CASH 38664 -> 144
$ws.Range("F2").Value = 144

# J2: REALIZED_PROFIT 35892 -> -240
$ws.Range("J2").Value = -240

# --- Row 3 ---
# B3: date text changes from 09/02/20 -> 09/04/20 (keep as text)
$ws.Range("B3").Value = "'09/04/20"
$ws.Range("B3").Style = "Normal"

# C3: BUY/SELL BUY -> SELL
$ws.Range("C3").Value = "SELL"

# D3: PRICE 23 -> 123
$ws.Range("D3").Value = 123

# F3: NET_EFFECT_TO_CASH -2806 -> 15006
$ws.Range("F3").Value = 15006

# G3: TOTAL_SHARES_HOLDING 122 -> 0
$ws.Range("G3").Value = 0

# H3: TICKER_TOTAL_VALUE 2806 -> 0
$ws.Range("H3").Value = 0

# I3: AVERAGE_PRICE 23 -> 0
$ws.Range("I3").Value = 0

# J3: REALIZED_PROFIT empty -> 11102 (numeric)
$ws.Range("J3").Value = 11102
